# Estadisticos Matutinos 15 Oct
# Fill in real Aprobados/Reprobados/Por_Apro/Por_Repro/Promedio values
# for rows 4-9 (Herrera Serrano Mayra Iliana, groups 1AM-1FM) on the
# "1er Parcial" and "3er Parcial" sheets.

$wb = $excel.ActiveWorkbook

$data = @{
    4 = @{ E = 25; F = 9;  G = 73.53; H = 26.47; I = 7.3; J = 9;  K = 26.47 }
    5 = @{ E = 22; F = 9;  G = 70.97; H = 29.03; I = 7.4; J = 9;  K = 29.03 }
    6 = @{ E = 29; F = 7;  G = 80.56; H = 19.44; I = 7.6; J = 7;  K = 19.44 }
    7 = @{ E = 32; F = 11; G = 74.42; H = 25.58; I = 7.5; J = 11; K = 25.58 }
    8 = @{ E = 33; F = 11; G = 75;    H = 25;    I = 7.4; J = 11; K = 25 }
    9 = @{ E = 23; F = 1;  G = 95.83; H = 4.17;  I = 7.3; J = 1;  K = 4.17 }
}

$sheetNames = @("1er Parcial", "3er Parcial")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $data.Keys) {
        $vals = $data[$row]
        $ws.Range("E$row").Value = $vals.E
        $ws.Range("F$row").Value = $vals.F
        $ws.Range("G$row").Value = $vals.G
        $ws.Range("H$row").Value = $vals.H
        $ws.Range("I$row").Value = $vals.I
        $ws.Range("J$row").Value = $vals.J
        $ws.Range("K$row").Value = $vals.K
    }
}
